# Enemies.xlsx update:
#  - extend enemy stat table (rows 18-41, columns B:F) with hpMax/damageMax/
#    defence/staminaMax/agressivity values that were previously blank
#  - fix the saved view state so opening the sheet no longer leaves it
#    scrolled to the bottom with row 33 selected (looked like the user had
#    stopped/waited there, when they'd actually been scrolling/"running"
#    through the list)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stat rows: row number -> (hpMax, damageMax, defence, staminaMax, agressivity)
$data = @{
    18 = @(4, 4, 0, 2, 7)     # Skeleton
    19 = @(20, 3, 3, 5, 1)    # Ent
    20 = @(6, 5, 1, 6, 4)     # Kobold
    21 = @(12, 6, 1, 9, 7)    # Hag
    22 = @(35, 10, 2, 12, 8)  # Wyvern
    23 = @(18, 8, 0, 12, 8)   # Griffon
    24 = @(30, 10, 1, 15, 7)  # Manticore
    25 = @(12, 4, 1, 10, 2)   # Centaur
    26 = @(8, 4, 0, 0, 6)     # Carnivorous Plant
    27 = @(12, 7, 2, 6, 9)    # Mad Knight
    28 = @(15, 5, 0, 10, 7)   # Vampire
    29 = @(13, 4, 1, 10, 7)   # Gargoyle
    30 = @(10, 3, 0, 12, 5)   # Cannibal
    31 = @(25, 8, 2, 10, 6)   # Basilisk
    32 = @(15, 10, 0, 15, 8)  # Acromantula
    33 = @(9, 5, 0, 20, 5)    # Succubus
    34 = @(8, 5, 1, 10, 7)    # Specter
    35 = @(4, 4, 0, 0, 6)     # Toxic Spore
    36 = @(10, 4, 0, 10, 6)   # Fungoid
    37 = @(15, 12, 2, 15, 8)  # Shadow Demon
    38 = @(10, 5, 1, 10, 8)   # Infernal Fiend
    39 = @(28, 8, 0, 10, 7)   # Naga
    40 = @(12, 8, 0, 15, 7)   # Sea Serpent
    41 = @(45, 14, 2, 15, 9)  # The Krakken
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
    }
}

# Fix the saved view/selection so the sheet no longer opens scrolled all the
# way down with A33 selected -- scroll/select near the top instead.
[void]$excel.Goto($ws.Range("A7"), $true)
[void]$ws.Range("F38").Select()
